$wb = $excel.ActiveWorkbook

# Sheets (by position, matching the workbook's tab order):
#   1 = F24 % cover
#   2 = F24 trees
#   3 = Historical Frequency
#   4 = Historical Rel. Frequency
#   5 = Historical Avg % Cover
$wsCover  = $wb.Worksheets.Item(1)
$wsFreq   = $wb.Worksheets.Item(3)
$wsRelFreq = $wb.Worksheets.Item(4)

# --- "F24 % cover": move the selection (no change of active tab here) ---
$wsCover.Range("H11").Select()

# --- "Historical Rel. Frequency": add the new "Relative Frequency" column (K) ---
# New historical relative-frequency values, computed as (count/30)*100,
# formatted to one decimal place.
$wsRelFreq.Range("K3").NumberFormat = "0.0"
$wsRelFreq.Range("K3").Formula = "=(5/30)*100"

$wsRelFreq.Range("K4").NumberFormat = "0.0"
$wsRelFreq.Range("K4").Formula = "=(7/30)*100"

$wsRelFreq.Range("K5").NumberFormat = "0.0"
$wsRelFreq.Range("K5").Value = 0

$wsRelFreq.Range("K6").NumberFormat = "0.0"
$wsRelFreq.Range("K6").Value = 0

$wsRelFreq.Range("K7").NumberFormat = "0.0"
$wsRelFreq.Range("K7").Formula = "=(5/30)*100"

$wsRelFreq.Range("K8").NumberFormat = "0.0"
$wsRelFreq.Range("K8").Formula = "=(10/30)*100"

$wsRelFreq.Range("K9").NumberFormat = "0.0"
$wsRelFreq.Range("K9").Formula = "=(3/30)*100"

# --- "Historical Frequency": move the selection ---
$wsFreq.Activate()
$wsFreq.Range("G18").Select()

# --- Make "Historical Rel. Frequency" the active tab/selection last, so it
#     ends up the active sheet in the saved workbook ---
$wsRelFreq.Activate()
$wsRelFreq.Range("F17").Select()
